$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("first_eval")

# New label for the added row (row 10) - copy formatting from row 9's label cell, then set the text
$ws.Range("A9").Copy($ws.Range("A10"))
$ws.Range("A10").Value = "Q8"

# Updated / corrected numeric values for rows 2-9 (columns B:F)
$values = @{
    2  = @(-0.3989960852331779, 1.622347057549135, 3.410040606493075, 1.846629526053636, 1.871071232177386)
    3  = @(-0.07368367304384635, 2.01870334081424, 6.614665517794586, 2.571899204439121, 2.675818738311348)
    4  = @(-0.3412105586050318, 1.726870111009629, 4.898523348216152, 2.213260795346123, 2.284039174058081)
    5  = @(-0.4693883225059389, 2.197338286322757, 7.681313688700947, 2.771518300264486, 2.864801410313742)
    6  = @(-0.7369237006858437, 1.317817770111458, 2.080685236058409, 1.442458053483154, 1.30708696277294)
    7  = @(-0.2355834118156171, 1.992865734783982, 4.408914335405655, 2.099741492518937, 2.213050302279264)
    8  = @(0.04677381843788716, 1.625585584899291, 3.112956264121649, 1.764357181559802, 1.932077164306971)
    9  = @(0.8102998313761515, 0.9132997309681818, 1.270800624281741, 1.127297930576359, 0.9598553074970639)
}

foreach ($row in $values.Keys) {
    $rowValues = $values[$row]
    $ws.Range("B$row").Value = $rowValues[0]
    $ws.Range("C$row").Value = $rowValues[1]
    $ws.Range("D$row").Value = $rowValues[2]
    $ws.Range("E$row").Value = $rowValues[3]
    $ws.Range("F$row").Value = $rowValues[4]
}

# Updated rank values in column G for rows 7 and 8
$ws.Range("G7").Value = 9
$ws.Range("G8").Value = 6

# New row 10 values (B:E populated, F left empty, G rank = 1)
$ws.Range("B10").Value = 0.1034734828819666
$ws.Range("C10").Value = 0.1034734828819666
$ws.Range("D10").Value = 0.01070676165972463
$ws.Range("E10").Value = 0.1034734828819666
$ws.Range("G10").Value = 1

# Row 9's previous rank value of 1 moves to row 10; row 9 is now rank 3
$ws.Range("G9").Value = 3
